# Increment the study-number labels in column A (rows 20-66) by 1.
# These cells hold the running study index that precedes each author's
# name; a new study was inserted earlier in the list, so every following
# index shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 20; $row -le 66; $row++) {
    $cell = $ws.Range("A" + $row)
    $current = $cell.Value2
    if ($null -ne $current -and $current -ne "") {
        $newValue = [int]$current + 1
        # Prefix with an apostrophe so Excel stores the result as text
        # (matching the original cell's text representation), then reset
        # the cell style back to Normal so no numeric/text formatting is
        # left behind on the cell.
        $cell.Value = "'" + [string]$newValue
        $cell.Style = "Normal"
    }
}
